# modified takscreenshot function in DSL
# The "TakeNativeScreenshot" calls used for VT200_0578 / VT200_0595 / VT200_0596
# are switched over to the regular "TakeScreenshot" helper, matching the
# other Steps cells in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

foreach ($addr in @("G4", "G6", "G7")) {
    $rng = $ws.Range($addr)
    $text = $rng.Value()
    $newText = $text.Replace("TakeNativeScreenshot(", "TakeScreenshot(")
    $rng.Value = $newText
}

# Reflect where the author was last working: scrolled down a couple of
# rows with G4 as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G4").Select()
